# Update NATMI TPM-derived ligand-receptor edge statistics for Angptl1-Tek
# across all sending/target cluster combinations (ECs, FAPs, MuSCs, Resolving-Mac).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Angptl1"
$ws.Range("C2").Value = "Tek"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.7697463333333333
$ws.Range("H2").Value = 2.309239
$ws.Range("I2").Value = 0.004764162270914757
$ws.Range("J2").Value = 0.004764162270914756
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 56.16842399999999
$ws.Range("N2").Value = 168.505272
$ws.Range("O2").Value = 0.8021403195141494
$ws.Range("P2").Value = 0.8021403195141494
$ws.Range("Q2").Value = 43.23543842311199
$ws.Range("R2").Value = 389.1189458080079
$ws.Range("S2").Value = 0.003821526646208819
$ws.Range("T2").Value = 0.003821526646208818

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Angptl1"
$ws.Range("C3").Value = "Tek"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.7697463333333333
$ws.Range("H3").Value = 2.309239
$ws.Range("I3").Value = 0.004764162270914757
$ws.Range("J3").Value = 0.004764162270914756
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.582537
$ws.Range("N3").Value = 31.747611
$ws.Range("O3").Value = 0.1511290331103168
$ws.Range("P3").Value = 0.1511290331103167
$ws.Range("Q3").Value = 8.145869053114334
$ws.Range("R3").Value = 73.31282147802899
$ws.Range("S3").Value = 0.000720003237583998
$ws.Range("T3").Value = 0.0007200032375839978

# Row 4: ECs -> MuSCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Angptl1"
$ws.Range("C4").Value = "Tek"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.7697463333333333
$ws.Range("H4").Value = 2.309239
$ws.Range("I4").Value = 0.004764162270914757
$ws.Range("J4").Value = 0.004764162270914756
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.262797333333333
$ws.Range("N4").Value = 9.788392
$ws.Range("O4").Value = 0.04659595390231912
$ws.Range("P4").Value = 0.04659595390231912
$ws.Range("Q4").Value = 2.511526283743111
$ws.Range("R4").Value = 22.603736553688
$ws.Range("S4").Value = 0.000221990685558712
$ws.Range("T4").Value = 0.0002219906855587119

# Row 5: ECs -> Resolving-Mac
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Angptl1"
$ws.Range("C5").Value = "Tek"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.7697463333333333
$ws.Range("H5").Value = 2.309239
$ws.Range("I5").Value = 0.004764162270914757
$ws.Range("J5").Value = 0.004764162270914756
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.009431666666666666
$ws.Range("N5").Value = 0.028295
$ws.Range("O5").Value = 0.0001346934732146117
$ws.Range("P5").Value = 0.0001346934732146117
$ws.Range("Q5").Value = 0.007259990833888889
$ws.Range("R5").Value = 0.065339917505
$ws.Range("S5").Value = 0.0000006417015632275203
$ws.Range("T5").Value = 0.0000006417015632275203

# Row 6: FAPs -> ECs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Angptl1"
$ws.Range("C6").Value = "Tek"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 159.049647
$ws.Range("H6").Value = 477.148941
$ws.Range("I6").Value = 0.9844000479461552
$ws.Range("J6").Value = 0.984400047946155
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 56.16842399999999
$ws.Range("N6").Value = 168.505272
$ws.Range("O6").Value = 0.8021403195141494
$ws.Range("P6").Value = 0.8021403195141494
$ws.Range("Q6").Value = 8933.568009746328
$ws.Range("R6").Value = 80402.11208771696
$ws.Range("S6").Value = 0.789626968989273
$ws.Range("T6").Value = 0.7896269689892728

# Row 7: FAPs -> FAPs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Angptl1"
$ws.Range("C7").Value = "Tek"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 159.049647
$ws.Range("H7").Value = 477.148941
$ws.Range("I7").Value = 0.9844000479461552
$ws.Range("J7").Value = 0.984400047946155
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 10.582537
$ws.Range("N7").Value = 31.747611
$ws.Range("O7").Value = 0.1511290331103168
$ws.Range("P7").Value = 0.1511290331103167
$ws.Range("Q7").Value = 1683.148774214439
$ws.Range("R7").Value = 15148.33896792995
$ws.Range("S7").Value = 0.1487714274398519
$ws.Range("T7").Value = 0.1487714274398518

# Row 8: FAPs -> MuSCs
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Angptl1"
$ws.Range("C8").Value = "Tek"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 159.049647
$ws.Range("H8").Value = 477.148941
$ws.Range("I8").Value = 0.9844000479461552
$ws.Range("J8").Value = 0.984400047946155
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.262797333333333
$ws.Range("N8").Value = 9.788392
$ws.Range("O8").Value = 0.04659595390231912
$ws.Range("P8").Value = 0.04659595390231912
$ws.Range("Q8").Value = 518.9467640992081
$ws.Range("R8").Value = 4670.520876892872
$ws.Range("S8").Value = 0.04586905925553979
$ws.Range("T8").Value = 0.04586905925553977

# Row 9: FAPs -> Resolving-Mac
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Angptl1"
$ws.Range("C9").Value = "Tek"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 159.049647
$ws.Range("H9").Value = 477.148941
$ws.Range("I9").Value = 0.9844000479461552
$ws.Range("J9").Value = 0.984400047946155
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.009431666666666666
$ws.Range("N9").Value = 0.028295
$ws.Range("O9").Value = 0.0001346934732146117
$ws.Range("P9").Value = 0.0001346934732146117
$ws.Range("Q9").Value = 1.500103253955
$ws.Range("R9").Value = 13.500929285595
$ws.Range("S9").Value = 0.0001325922614904979
$ws.Range("T9").Value = 0.0001325922614904979

# Row 10: MuSCs -> ECs
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Angptl1"
$ws.Range("C10").Value = "Tek"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.407514666666667
$ws.Range("H10").Value = 4.222544
$ws.Range("I10").Value = 0.008711478028942643
$ws.Range("J10").Value = 0.008711478028942642
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 56.16842399999999
$ws.Range("N10").Value = 168.505272
$ws.Range("O10").Value = 0.8021403195141494
$ws.Range("P10").Value = 0.8021403195141494
$ws.Range("Q10").Value = 79.05788058355199
$ws.Range("R10").Value = 711.5209252519679
$ws.Range("S10").Value = 0.006987827769576544
$ws.Range("T10").Value = 0.006987827769576543

# Row 11: MuSCs -> FAPs
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Angptl1"
$ws.Range("C11").Value = "Tek"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.407514666666667
$ws.Range("H11").Value = 4.222544
$ws.Range("I11").Value = 0.008711478028942643
$ws.Range("J11").Value = 0.008711478028942642
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 10.582537
$ws.Range("N11").Value = 31.747611
$ws.Range("O11").Value = 0.1511290331103168
$ws.Range("P11").Value = 0.1511290331103167
$ws.Range("Q11").Value = 14.89507603804267
$ws.Range("R11").Value = 134.055684342384
$ws.Range("S11").Value = 0.00131655725147587
$ws.Range("T11").Value = 0.001316557251475869

# Row 12: MuSCs -> MuSCs
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Angptl1"
$ws.Range("C12").Value = "Tek"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.407514666666667
$ws.Range("H12").Value = 4.222544
$ws.Range("I12").Value = 0.008711478028942643
$ws.Range("J12").Value = 0.008711478028942642
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 3.262797333333333
$ws.Range("N12").Value = 9.788392
$ws.Range("O12").Value = 0.04659595390231912
$ws.Range("P12").Value = 0.04659595390231912
$ws.Range("Q12").Value = 4.592435101027555
$ws.Range("R12").Value = 41.331915909248
$ws.Range("S12").Value = 0.0004059196286576773
$ws.Range("T12").Value = 0.0004059196286576772

# Row 13: MuSCs -> Resolving-Mac
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Angptl1"
$ws.Range("C13").Value = "Tek"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.407514666666667
$ws.Range("H13").Value = 4.222544
$ws.Range("I13").Value = 0.008711478028942643
$ws.Range("J13").Value = 0.008711478028942642
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.009431666666666666
$ws.Range("N13").Value = 0.028295
$ws.Range("O13").Value = 0.0001346934732146117
$ws.Range("P13").Value = 0.0001346934732146117
$ws.Range("Q13").Value = 0.01327520916444444
$ws.Range("R13").Value = 0.11947688248
$ws.Range("S13").Value = 0.000001173379232551064
$ws.Range("T13").Value = 0.000001173379232551064

# Row 14: Resolving-Mac -> ECs
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Angptl1"
$ws.Range("C14").Value = "Tek"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.3432253333333333
$ws.Range("H14").Value = 1.029676
$ws.Range("I14").Value = 0.002124311753987536
$ws.Range("J14").Value = 0.002124311753987535
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 56.16842399999999
$ws.Range("N14").Value = 168.505272
$ws.Range("O14").Value = 0.8021403195141494
$ws.Range("P14").Value = 0.8021403195141494
$ws.Range("Q14").Value = 19.278426050208
$ws.Range("R14").Value = 173.505834451872
$ws.Range("S14").Value = 0.001703996109091225
$ws.Range("T14").Value = 0.001703996109091225

# Row 15: Resolving-Mac -> FAPs
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Angptl1"
$ws.Range("C15").Value = "Tek"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.3432253333333333
$ws.Range("H15").Value = 1.029676
$ws.Range("I15").Value = 0.002124311753987536
$ws.Range("J15").Value = 0.002124311753987535
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 10.582537
$ws.Range("N15").Value = 31.747611
$ws.Range("O15").Value = 0.1511290331103168
$ws.Range("P15").Value = 0.1511290331103167
$ws.Range("Q15").Value = 3.632194789337333
$ws.Range("R15").Value = 32.689753104036
$ws.Range("S15").Value = 0.0003210451814050174
$ws.Range("T15").Value = 0.0003210451814050173

# Row 16: Resolving-Mac -> MuSCs
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Angptl1"
$ws.Range("C16").Value = "Tek"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.3432253333333333
$ws.Range("H16").Value = 1.029676
$ws.Range("I16").Value = 0.002124311753987536
$ws.Range("J16").Value = 0.002124311753987535
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 3.262797333333333
$ws.Range("N16").Value = 9.788392
$ws.Range("O16").Value = 0.04659595390231912
$ws.Range("P16").Value = 0.04659595390231912
$ws.Range("Q16").Value = 1.119874702332444
$ws.Range("R16").Value = 10.078872320992
$ws.Range("S16").Value = 0.00009898433256295791
$ws.Range("T16").Value = 0.00009898433256295788

# Row 17: Resolving-Mac -> Resolving-Mac
$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Angptl1"
$ws.Range("C17").Value = "Tek"
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.3432253333333333
$ws.Range("H17").Value = 1.029676
$ws.Range("I17").Value = 0.002124311753987536
$ws.Range("J17").Value = 0.002124311753987535
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.009431666666666666
$ws.Range("N17").Value = 0.028295
$ws.Range("O17").Value = 0.0001346934732146117
$ws.Range("P17").Value = 0.0001346934732146117
$ws.Range("Q17").Value = 0.003237186935555555
$ws.Range("R17").Value = 0.02913468242
$ws.Range("S17").Value = 0.0000002861309283352049
$ws.Range("T17").Value = 0.0000002861309283352049

